$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round all numeric data values in B2:E13 to the nearest integer,
# matching the "write to disk as integer data" behavior described
# in the commit message.
$rng = $ws.Range("B2:E13")
foreach ($cell in $rng.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = [Math]::Floor([double]$val + 0.5)
    }
}
